$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.046534378652948
$ws.Range("D2").Value = 1.049846645543807
$ws.Range("E2").Value = 1.044084260354381
$ws.Range("F2").Value = 1.05901497523259
$ws.Range("I2").Value = 1.039146763037203
$ws.Range("J2").Value = 1.051588260669239
$ws.Range("K2").Value = 1.052602567067423
$ws.Range("L2").Value = 1.046856306123832
$ws.Range("M2").Value = 1.06174563383294
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.047979743052009
$ws.Range("D3").Value = 1.050927493606549
$ws.Range("E3").Value = 1.045330155527064
$ws.Range("F3").Value = 1.060243561956825
$ws.Range("I3").Value = 1.039455476554014
$ws.Range("J3").Value = 1.052679290317649
$ws.Range("K3").Value = 1.053495050327574
$ws.Range("L3").Value = 1.047912232682433
$ws.Range("M3").Value = 1.062787317718793
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.048914098613176
$ws.Range("D4").Value = 1.051625910123484
$ws.Range("E4").Value = 1.046135779253045
$ws.Range("F4").Value = 1.061037832366607
$ws.Range("I4").Value = 1.039653464321051
$ws.Range("J4").Value = 1.053383952083915
$ws.Range("K4").Value = 1.054070990328696
$ws.Range("L4").Value = 1.048594378973953
$ws.Range("M4").Value = 1.063460085431428
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.049306693708842
$ws.Range("D5").Value = 1.051919296636344
$ws.Range("E5").Value = 1.046474334767035
$ws.Range("F5").Value = 1.061371577847241
$ws.Range("I5").Value = 1.039736275678081
$ws.Range("J5").Value = 1.053679882693965
$ws.Range("K5").Value = 1.054312745855686
$ws.Range("L5").Value = 1.048880891362997
$ws.Range("M5").Value = 1.063742615850366
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.049372600117029
$ws.Range("D6").Value = 1.051968544254861
$ws.Range("E6").Value = 1.046531172292673
$ws.Range("F6").Value = 1.061427605562036
$ws.Range("I6").Value = 1.039750155315725
$ws.Range("J6").Value = 1.053729552723957
$ws.Range("K6").Value = 1.054353316081253
$ws.Range("L6").Value = 1.048928982797671
$ws.Range("M6").Value = 1.063790036386843
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.048919345298427
$ws.Range("D7").Value = 1.051629831259656
$ws.Range("E7").Value = 1.0461403035493
$ws.Range("F7").Value = 1.061042292537146
$ws.Range("I7").Value = 1.039654572509941
$ws.Range("J7").Value = 1.05338790753345
$ws.Range("K7").Value = 1.054074222126715
$ws.Range("L7").Value = 1.048598208391276
$ws.Range("M7").Value = 1.063463861795918
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047023034317174
$ws.Range("D8").Value = 1.050212124775588
$ws.Range("E8").Value = 1.044505433456517
$ws.Range("F8").Value = 1.059430330468118
$ws.Range("I8").Value = 1.0392514615179
$ws.Range("J8").Value = 1.051957252233892
$ws.Range("K8").Value = 1.052904509991854
$ws.Range("L8").Value = 1.047213393166275
$ws.Range("M8").Value = 1.062097940997232
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043674413369201
$ws.Range("D9").Value = 1.047706418553496
$ws.Range("E9").Value = 1.041620168270327
$ws.Range("F9").Value = 1.056584270431713
$ws.Range("I9").Value = 1.038527515737629
$ws.Range("J9").Value = 1.049426071465251
$ws.Range("K9").Value = 1.050831275495357
$ws.Range("L9").Value = 1.044764520677927
$ws.Range("M9").Value = 1.059681130096646
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041436891292179
$ws.Range("D10").Value = 1.0460306896912
$ws.Range("E10").Value = 1.039693444049454
$ws.Range("F10").Value = 1.054682922033517
$ws.Range("I10").Value = 1.038035657956309
$ws.Range("J10").Value = 1.047731530533521
$ws.Range("K10").Value = 1.049440828488295
$ws.Range("L10").Value = 1.043125898079315
$ws.Range("M10").Value = 1.058063070409045
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040466729805599
$ws.Range("D11").Value = 1.045303788674874
$ws.Range("E11").Value = 1.038858333535813
$ws.Range("F11").Value = 1.053858623424878
$ws.Range("I11").Value = 1.037820472467436
$ws.Range("J11").Value = 1.046996040513565
$ws.Range("K11").Value = 1.048836740780723
$ws.Range("L11").Value = 1.042414873016589
$ws.Range("M11").Value = 1.057360760616898
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040106166561445
$ws.Range("D12").Value = 1.045033586103834
$ws.Range("E12").Value = 1.038548007328177
$ws.Range("F12").Value = 1.053552287165508
$ws.Range("I12").Value = 1.037740209762215
$ws.Range("J12").Value = 1.046722580212254
$ws.Range("K12").Value = 1.048612049239862
$ws.Range("L12").Value = 1.042150538254103
$ws.Range("M12").Value = 1.057099635072883
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040183517855242
$ws.Range("D13").Value = 1.045091554549124
$ws.Range("E13").Value = 1.038614579267958
$ws.Range("F13").Value = 1.053618004421884
$ws.Range("I13").Value = 1.037757441489192
$ws.Range("J13").Value = 1.046781250514203
$ws.Range("K13").Value = 1.048660260266839
$ws.Range("L13").Value = 1.042207249364953
$ws.Range("M13").Value = 1.05715565905843
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040436929697101
$ws.Range("D14").Value = 1.045281457732007
$ws.Range("E14").Value = 1.038832684533996
$ws.Range("F14").Value = 1.053833304763222
$ws.Range("I14").Value = 1.037813844733953
$ws.Range("J14").Value = 1.046973441665094
$ws.Range("K14").Value = 1.048818173981477
$ws.Range("L14").Value = 1.042393027707654
$ws.Range("M14").Value = 1.05733918115276
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040593038115611
$ws.Range("D15").Value = 1.045398436840073
$ws.Range("E15").Value = 1.038967049150098
$ws.Range("F15").Value = 1.05396593772847
$ws.Range("I15").Value = 1.03784855243187
$ws.Range("J15").Value = 1.047091821511775
$ws.Range("K15").Value = 1.048915429110353
$ws.Range("L15").Value = 1.042507461478068
$ws.Range("M15").Value = 1.05745222107045
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041501249751965
$ws.Range("D16").Value = 1.046078904004776
$ws.Range("E16").Value = 1.039748849763235
$ws.Range("F16").Value = 1.054737606520047
$ws.Range("I16").Value = 1.038049892461604
$ws.Range("J16").Value = 1.047780305464063
$ws.Range("K16").Value = 1.049480877027884
$ws.Range("L16").Value = 1.043173054686087
$ws.Range("M16").Value = 1.058109644665168
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042070594163166
$ws.Range("D17").Value = 1.04650539242074
$ws.Range("E17").Value = 1.040239028017601
$ws.Range("F17").Value = 1.05522138249037
$ws.Range("I17").Value = 1.038175595566847
$ws.Range("J17").Value = 1.04821170336474
$ws.Range("K17").Value = 1.049835025467019
$ws.Range("L17").Value = 1.043590161501385
$ws.Range("M17").Value = 1.058521576624886
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042402558190089
$ws.Range("D18").Value = 1.046754030798145
$ws.Range("E18").Value = 1.040524861312912
$ws.Range("F18").Value = 1.055503464436186
$ws.Range("I18").Value = 1.038248703142956
$ws.Range("J18").Value = 1.048463162647868
$ws.Range("K18").Value = 1.050041400098525
$ws.Range("L18").Value = 1.043833309324555
$ws.Range("M18").Value = 1.058761687823444
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042515728358449
$ws.Range("D19").Value = 1.046838788959134
$ws.Range("E19").Value = 1.04062230976046
$ws.Range("F19").Value = 1.055599630874961
$ws.Range("I19").Value = 1.038273594868518
$ws.Range("J19").Value = 1.048548875443986
$ws.Range("K19").Value = 1.050111735694027
$ws.Range("L19").Value = 1.043916192237944
$ws.Range("M19").Value = 1.058843532226527
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042009521868226
$ws.Range("D20").Value = 1.046459647218414
$ws.Range("E20").Value = 1.04018644477659
$ws.Range("F20").Value = 1.055169487910162
$ws.Range("I20").Value = 1.038162130847615
$ws.Range("J20").Value = 1.048165435845011
$ws.Range("K20").Value = 1.049797048801363
$ws.Range("L20").Value = 1.043545424765348
$ws.Range("M20").Value = 1.058477396984073
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.04036231185887
$ws.Range("D21").Value = 1.045225541507779
$ws.Range("E21").Value = 1.03876846156886
$ws.Range("F21").Value = 1.053769908495341
$ws.Range("I21").Value = 1.037797244595428
$ws.Range("J21").Value = 1.046916853558079
$ws.Range("K21").Value = 1.048771680818642
$ws.Range("L21").Value = 1.042338326956779
$ws.Range("M21").Value = 1.057285145581784
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039325473021557
$ws.Range("D22").Value = 1.044448455960253
$ws.Range("E22").Value = 1.037876170585039
$ws.Range("F22").Value = 1.052889038028634
$ws.Range("I22").Value = 1.03756589737264
$ws.Range("J22").Value = 1.046130275707091
$ws.Range("K22").Value = 1.048125216262427
$ws.Range("L22").Value = 1.04157805300702
$ws.Range("M22").Value = 1.056534044191862
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.039875234332619
$ws.Range("D23").Value = 1.044860514453632
$ws.Range("E23").Value = 1.038349263509604
$ws.Range("F23").Value = 1.053356090881004
$ws.Range("I23").Value = 1.037688722250621
$ws.Range("J23").Value = 1.046547403513509
$ws.Range("K23").Value = 1.048468088860798
$ws.Range("L23").Value = 1.041981215568443
$ws.Range("M23").Value = 1.056932359452956
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.042037118195877
$ws.Range("D24").Value = 1.046480317893179
$ws.Range("E24").Value = 1.040210205125457
$ws.Range("F24").Value = 1.055192937135426
$ws.Range("I24").Value = 1.038168215632697
$ws.Range("J24").Value = 1.048186342664471
$ws.Range("K24").Value = 1.049814209422129
$ws.Range("L24").Value = 1.043565639814495
$ws.Range("M24").Value = 1.058497360362187
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.04454098680985
$ws.Range("D25").Value = 1.048355115852485
$ws.Range("E25").Value = 1.042366626893345
$ws.Range("F25").Value = 1.057320729181082
$ws.Range("I25").Value = 1.03871629392121
$ws.Range("J25").Value = 1.050081673801865
$ws.Range("K25").Value = 1.051368703779197
$ws.Range("L25").Value = 1.04539866113666
$ws.Range("M25").Value = 1.060307126169891
